# Update "productos" catalog: refresh prices, unify talla (size) range,
# rename a few product rows, and highlight the updated price cells in yellow.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: a handful of product names were renamed (rows 12-16) ---
$ws.Range("A12").Value = "UNIVERSAL"
$ws.Range("A13").Value = "NEO"
$ws.Range("A14").Value = "BARRIO"
$ws.Range("A15").Value = "VALO"
$ws.Range("A16").Value = "ONE"

# --- Column B: new prices; row 4 (PRO MAX FUTSAL) keeps its price ---
$ws.Range("B2").Value = 43000
$ws.Range("B3").Value = 43000
$ws.Range("B5").Value = 38000
$ws.Range("B6").Value = 37000
$ws.Range("B7").Value = 45000
$ws.Range("B8").Value = 38000
$ws.Range("B9").Value = 38000
$ws.Range("B10").Value = 43000
$ws.Range("B11").Value = 37000
$ws.Range("B12").Value = 37000
$ws.Range("B13").Value = 37000
$ws.Range("B14").Value = 37000
$ws.Range("B15").Value = 37000
$ws.Range("B16").Value = 37000

# Highlight every updated price cell (all but B4) with a solid yellow fill.
$ws.Range("B2").Interior.Color = 65535
$ws.Range("B3").Interior.Color = 65535
$ws.Range("B5:B16").Interior.Color = 65535

# --- Column C: every row now shares the same talla (size) range ---
$ws.Range("C2:C16").Value = "39,40,41,42,43"

# Selection now spans the whole tallas column, matching the saved view state.
$ws.Range("C2:C16").Select() | Out-Null
